# "run 2019 final report" - append new survey response rows to the research metrics log,
# update one quantity value, and move the active selection to the new last entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item("Table1")

# --- 1. Update existing Quantity value in row 250 (13 -> 16) ---
$ws.Range("F250").Value2 = 16

# --- 2. Append 9 new rows of data to the table (rows 253-261) ---
$newData = @(
    @{ A = 43843.74900462963;  B = "ahappel@sheddaquarium.org"; C = "Citizen/stakeholder engagement"; D = 43838; E = "https://www.reddit.com/r/MicroFishing/comments/elty76/rmicrofishing_shout_out_in_research_paper/"; F = 50;  G = "Posted to the Reddit Microfishing page that we spoke about them in a research paper. It received 50 upvotes, which means ATLEAST 50 people read my post, but likely way more did as the community numbers > 5,000" },
    @{ A = 43858.879895833335; B = "ahappel@sheddaquarium.org"; C = "Citizen/stakeholder engagement"; D = 43852; E = $null; F = 15; G = "Lecture and Lab for vertebrate ecology course at Loyola." },
    @{ A = 43859.769270833334; B = "kinoue@sheddaquarium.org";  C = "Publication"; D = 43852; E = "Inoue K, Pohl AL, Makiri S, Lang BK, Berg DJ. (2020) Use of species delimitation approach to assess biodiversity in freshwater planarians (Platyhelminthes: Tricladida) from desert springs. Aquatic Conservation: Marine and Freshwater Ecosystems, DOI: 10.1002/aqc.3273"; F = $null; G = $null },
    @{ A = 43861.711574074077; B = "ahappel@sheddaquarium.org"; C = "Citizen/stakeholder engagement"; D = 43860; E = $null; F = 125; G = "Presented at The Night of Ideas. >5,000 attended event, and estimated 125 came to my talk." },
    @{ A = 43861.778101851851; B = "skessel@sheddaquarium.org"; C = "Citizen/stakeholder engagement"; D = 43790; E = $null; F = 50; G = "Spoke about our research programs at the 'Waves of Gratitude' event for the auxiliary board" },
    @{ A = 43861.778796296298; B = "skessel@sheddaquarium.org"; C = "Field research"; D = 43743; E = $null; F = 7; G = "Coral research trip" },
    @{ A = 43861.779282407406; B = "AKOUGH@sheddaquarium.org";  C = "Field research"; D = 43749; E = $null; F = 7; G = "Coral trip" },
    @{ A = 43861.779849537037; B = "AKOUGH@sheddaquarium.org";  C = "Media opportunity"; D = 43861; E = $null; F = $null; G = "Interview for ScienceWorld - a Scholastic magazine for highschoolers" },
    @{ A = 43861.918692129628; B = "kmurchie@sheddaquarium.org"; C = "Publication"; D = 43826; E = "Happel A., K.J. Murchie, P. W. Willink, and C.R. Knapp. In Press. Great Lakes Fish Finder App; a tool for biologists, managers and education practitioners. Journal of Great Lakes Research. XX:XX-XX. https://doi.org/10.1016/j.jglr.2019.12.002 "; F = $null; G = $null }
)

foreach ($row in $newData) {
    $listRow = $tbl.ListRows.Add()
    $rowIndex = $listRow.Range.Row

    # copy number formats (date / datetime styles) from the row directly above
    $ws.Range("A" + ($rowIndex - 1) + ":G" + ($rowIndex - 1)).Copy() | Out-Null
    $ws.Range("A" + $rowIndex + ":G" + $rowIndex).PasteSpecial(-4122) | Out-Null

    $ws.Range("A" + $rowIndex).Value2 = $row.A
    $ws.Range("B" + $rowIndex).Value2 = $row.B
    $ws.Range("C" + $rowIndex).Value2 = $row.C
    $ws.Range("D" + $rowIndex).Value2 = $row.D
    if ($row.E -ne $null) { $ws.Range("E" + $rowIndex).Value2 = $row.E } else { $ws.Range("E" + $rowIndex).Clear() }
    if ($row.F -ne $null) { $ws.Range("F" + $rowIndex).Value2 = $row.F } else { $ws.Range("F" + $rowIndex).Clear() }
    if ($row.G -ne $null) { $ws.Range("G" + $rowIndex).Value2 = $row.G } else { $ws.Range("G" + $rowIndex).Clear() }
}

$excel.CutCopyMode = 0

# --- 3. Move the active selection/view to the new bottom of the data (row 251 in view) ---
$ws.Range("F251").Select()
